$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96, shifting existing rows 96-173 down to 97-174
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new data record
$ws.Cells.Item(96, 1).Value = 3
$ws.Cells.Item(96, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(96, 3).Value = "Coquimbo"
$ws.Cells.Item(96, 4).Value = 44669
$ws.Cells.Item(96, 5).Value = 5
$ws.Cells.Item(96, 6).Value = 100112052
$ws.Cells.Item(96, 7).Value = "Albahaca"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 50
$ws.Cells.Item(96, 11).Value = 4500
$ws.Cells.Item(96, 12).Value = 4500
$ws.Cells.Item(96, 13).Value = 4500
$ws.Cells.Item(96, 14).Value = "$/docena de matas"
$ws.Cells.Item(96, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(96, 16).Value = 750
$ws.Cells.Item(96, 17).Value = 6
$ws.Cells.Item(96, 18).Value = "Hortaliza"
